$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "sample_rate"
$ws.Range("B9").Value = "uint16"
$ws.Range("C9").Value = "sampling rate of the recording"

$ws.Range("A10").Value = "epoch_total"
$ws.Range("B10").Value = "uint16"
$ws.Range("C10").Value = "the total number of the epochs"

[void]$ws.Range("H17").Select()
